$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Add paragraph border (top/left/bottom/right, space=5) to the first paragraph.
$b = $p1.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Replace the placeholder id text and drop the trailing space run.
$r = $p1.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "**ID__AFFARS_SUBPART_5303_10__ID**"
